$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 61086
$ws.Range("B2").Value = "Breno Cirino"
$ws.Range("C2").Value = "Operacoes"
$ws.Range("D2").Value = "Viagem de negocios"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 45097
$ws.Range("G2").Value = 9072.370000000001

# Row 3
$ws.Range("A3").Value = 9916
$ws.Range("B3").Value = "Guilherme Sousa"
$ws.Range("C3").Value = "Atendimento ao Cliente"
$ws.Range("D3").Value = "Consulta medica"
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 45103
$ws.Range("G3").Value = 4531.41

# Row 4
$ws.Range("A4").Value = 4248
$ws.Range("B4").Value = "Otto Pinto"
$ws.Range("C4").Value = "Engenharia"
$ws.Range("D4").Value = "Consulta medica"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 45088
$ws.Range("G4").Value = 9851.68

# Row 5
$ws.Range("A5").Value = 80295
$ws.Range("B5").Value = "Maria Luísa Macedo"
$ws.Range("C5").Value = "Operacoes"
$ws.Range("D5").Value = "Consulta medica"
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 45084
$ws.Range("G5").Value = 7867.06

# Row 6
$ws.Range("A6").Value = 27956
$ws.Range("B6").Value = "Isaque Marques"
$ws.Range("C6").Value = "P&D"
$ws.Range("D6").Value = "Viagem de negocios"
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 45090
$ws.Range("G6").Value = 3659.94

# Row 7
$ws.Range("A7").Value = 47
$ws.Range("B7").Value = "Olivia Macedo"
$ws.Range("C7").Value = "Atendimento ao Cliente"
$ws.Range("D7").Value = "Outros"
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 45106
$ws.Range("G7").Value = 9903.530000000001

# Row 8
$ws.Range("A8").Value = 22937
$ws.Range("B8").Value = "Luana Fernandes"
$ws.Range("C8").Value = "Operacoes"
$ws.Range("D8").Value = "Problemas pessoais"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 45091
$ws.Range("G8").Value = 3630.26

# Row 9
$ws.Range("A9").Value = 42081
$ws.Range("B9").Value = "Gael Leão"
$ws.Range("C9").Value = "TI"
$ws.Range("D9").Value = "Outros"
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 45090
$ws.Range("G9").Value = 3198.99

# Row 10
$ws.Range("A10").Value = 48812
$ws.Range("B10").Value = "Nathan Barbosa"
$ws.Range("C10").Value = "Juridico"
$ws.Range("D10").Value = "Viagem de negocios"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 45086
$ws.Range("G10").Value = 5212.71

# Row 11
$ws.Range("A11").Value = 78853
$ws.Range("B11").Value = "Fernando Monteiro"
$ws.Range("C11").Value = "P&D"
$ws.Range("D11").Value = "Outros"
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 45086
$ws.Range("G11").Value = 5205.82
